$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Revert "Speaker ner: added entities, widgets, binder":
# remove the last row (text4 / Chinese state media article) and
# rename the header from "text_name" to "text_id"
$ws.Rows.Item(5).Delete()
$ws.Range("A1").Value = "text_id"
$null = $ws.Range("B2").Select()
